$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): F4, F5, F6
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F4").Value = 1701
$ws1.Range("F5").Value = 764
$ws1.Range("F6").Value = 190

# Sheet "全部类型" (sheet4): F4, F6, F7
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F4").Value = 1701
$ws4.Range("F6").Value = 764
$ws4.Range("F7").Value = 190
